$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H40").Value = 1356.1333
$ws_ALC.Range("J40").Value = 2003
$ws_ALC.Range("L40").Value = 2003
$ws_ALC.Range("N40").Value = -2353
$ws_ALC.Range("H43").Value = 733.46155
$ws_ALC.Range("I43").Value = 596.3333
$ws_ALC.Range("J43").Value = 774.6
$ws_ALC.Range("K43").Value = 596.3333
$ws_ALC.Range("L43").Value = 774.6
$ws_ALC.Range("M43").Value = -527.3333
$ws_ALC.Range("N43").Value = -912.6
$ws_ALC.Range("H62").Value = 3607.2104
$ws_ALC.Range("I62").Value = 3163.125
$ws_ALC.Range("J62").Value = 3930.182
$ws_ALC.Range("K62").Value = 3163.125
$ws_ALC.Range("L62").Value = 3930.182
$ws_ALC.Range("M62").Value = -2539.125
$ws_ALC.Range("N62").Value = -5178.182
$ws_ALC.Range("H65").Value = 3607.2104
$ws_ALC.Range("I65").Value = 3163.125
$ws_ALC.Range("J65").Value = 3930.182
$ws_ALC.Range("K65").Value = 15815.625
$ws_ALC.Range("L65").Value = 19650.91
$ws_ALC.Range("M65").Value = -12695.625
$ws_ALC.Range("N65").Value = -25890.91
$ws_ALC.Range("H111").Value = 7962.6665
$ws_ALC.Range("J111").Value = 5000
$ws_ALC.Range("L111").Value = 15000
$ws_ALC.Range("N111").Value = -21134
$ws_ALC.Range("H116").Value = 8623601
$ws_ALC.Range("I116").Value = 20834380
$ws_ALC.Range("J116").Value = 4226.647
$ws_ALC.Range("K116").Value = 20834380
$ws_ALC.Range("L116").Value = 4226.647
$ws_ALC.Range("M116").Value = -20830938
$ws_ALC.Range("N116").Value = -11110.647
$ws_ALC.Range("H127").Value = 839.4
$ws_ALC.Range("I127").Value = 381.0909
$ws_ALC.Range("K127").Value = 1143.2727
$ws_ALC.Range("M127").Value = 3816.7273
$ws_ALC.Range("H129").Value = 127506.47
$ws_ALC.Range("J129").Value = 143867.7
$ws_ALC.Range("L129").Value = 431603.1
$ws_ALC.Range("N129").Value = -441603.1
$ws_ALC.Range("H138").Value = 2739.4614
$ws_ALC.Range("I138").Value = 2227
$ws_ALC.Range("J138").Value = 2876.9512
$ws_ALC.Range("K138").Value = 6681
$ws_ALC.Range("L138").Value = 8630.8536
$ws_ALC.Range("M138").Value = -1541
$ws_ALC.Range("N138").Value = -18910.8536

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H63").Value = 1839888.4
$ws_ARM.Range("I63").Value = 1759.9333
$ws_ARM.Range("J63").Value = 15625852
$ws_ARM.Range("K63").Value = 1759.9333
$ws_ARM.Range("L63").Value = 15625852
$ws_ARM.Range("M63").Value = -1073.9333
$ws_ARM.Range("N63").Value = -15627224
$ws_ARM.Range("H66").Value = 1839888.4
$ws_ARM.Range("I66").Value = 1759.9333
$ws_ARM.Range("J66").Value = 15625852
$ws_ARM.Range("K66").Value = 8799.666499999999
$ws_ARM.Range("L66").Value = 78129260
$ws_ARM.Range("M66").Value = -5367.666499999999
$ws_ARM.Range("N66").Value = -78136124
$ws_ARM.Range("H102").Value = 1232.6
$ws_ARM.Range("I102").Value = 1165.75
$ws_ARM.Range("K102").Value = 1165.75
$ws_ARM.Range("M102").Value = 456.25
$ws_ARM.Range("H132").Value = 14719.049
$ws_ARM.Range("I132").Value = 2474.3103
$ws_ARM.Range("J132").Value = 44310.5
$ws_ARM.Range("K132").Value = 7422.9309
$ws_ARM.Range("L132").Value = 132931.5
$ws_ARM.Range("M132").Value = -4892.9309
$ws_ARM.Range("N132").Value = -137991.5

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H94").Value = 733.93335
$ws_BSM.Range("I94").Value = 608.17645
$ws_BSM.Range("J94").Value = 898.38464
$ws_BSM.Range("K94").Value = 608.17645
$ws_BSM.Range("L94").Value = 898.38464
$ws_BSM.Range("M94").Value = -157.17645
$ws_BSM.Range("N94").Value = -1800.38464
$ws_BSM.Range("H105").Value = 1001843.4
$ws_BSM.Range("I105").Value = 1570.909
$ws_BSM.Range("J105").Value = 1787771.8
$ws_BSM.Range("K105").Value = 1570.909
$ws_BSM.Range("L105").Value = 1787771.8
$ws_BSM.Range("M105").Value = 176.0909999999999
$ws_BSM.Range("N105").Value = -1791265.8
$ws_BSM.Range("H132").Value = 0
$ws_BSM.Range("J132").Value = 0
$ws_BSM.Range("N132").ClearContents()

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H56").Value = 12134.286
$ws_CRP.Range("I56").Value = 2646.6667
$ws_CRP.Range("J56").Value = 19250
$ws_CRP.Range("K56").Value = 2646.6667
$ws_CRP.Range("L56").Value = 19250
$ws_CRP.Range("M56").Value = -1801.6667
$ws_CRP.Range("N56").Value = -20940
$ws_CRP.Range("H58").Value = 17124.156
$ws_CRP.Range("I58").Value = 1515.6842
$ws_CRP.Range("J58").Value = 39936.54
$ws_CRP.Range("K58").Value = 1515.6842
$ws_CRP.Range("L58").Value = 39936.54
$ws_CRP.Range("M58").Value = -1312.6842
$ws_CRP.Range("N58").Value = -40342.54
$ws_CRP.Range("H105").Value = 1071.4445
$ws_CRP.Range("I105").Value = 971.6667
$ws_CRP.Range("K105").Value = 971.6667
$ws_CRP.Range("M105").Value = 775.3333
$ws_CRP.Range("H134").Value = 1490.0588
$ws_CRP.Range("I134").Value = 911.5
$ws_CRP.Range("J134").Value = 2878.6
$ws_CRP.Range("K134").Value = 2734.5
$ws_CRP.Range("L134").Value = 8635.799999999999
$ws_CRP.Range("M134").Value = -199.5
$ws_CRP.Range("N134").Value = -13705.8
$ws_CRP.Range("H136").Value = 17124.156
$ws_CRP.Range("I136").Value = 1515.6842
$ws_CRP.Range("J136").Value = 39936.54
$ws_CRP.Range("K136").Value = 4547.0526
$ws_CRP.Range("L136").Value = 119809.62
$ws_CRP.Range("M136").Value = -1997.0526
$ws_CRP.Range("N136").Value = -124909.62

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H5").Value = 916.7451
$ws_CUL.Range("J5").Value = 1862.0714
$ws_CUL.Range("L5").Value = 5586.2142
$ws_CUL.Range("N5").Value = -5810.2142
$ws_CUL.Range("H12").Value = 115.77778
$ws_CUL.Range("I12").Value = 20
$ws_CUL.Range("J12").Value = 127.75
$ws_CUL.Range("K12").Value = 60
$ws_CUL.Range("L12").Value = 383.25
$ws_CUL.Range("M12").Value = 113
$ws_CUL.Range("N12").Value = -729.25
$ws_CUL.Range("H107").Value = 3751.2144
$ws_CUL.Range("I107").Value = 5178.75
$ws_CUL.Range("J107").Value = 182.375
$ws_CUL.Range("K107").Value = 15536.25
$ws_CUL.Range("L107").Value = 547.125
$ws_CUL.Range("M107").Value = -13616.25
$ws_CUL.Range("N107").Value = -4387.125
$ws_CUL.Range("H131").Value = 728.8099999999999
$ws_CUL.Range("J131").Value = 735.165
$ws_CUL.Range("L131").Value = 2205.495
$ws_CUL.Range("N131").Value = -12285.495
$ws_CUL.Range("H135").Value = 916.7451
$ws_CUL.Range("J135").Value = 1862.0714
$ws_CUL.Range("L135").Value = 16758.6426
$ws_CUL.Range("N135").Value = -21828.6426

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H80").Value = 3602.2917
$ws_GSM.Range("J80").Value = 3792.8572
$ws_GSM.Range("L80").Value = 3792.8572
$ws_GSM.Range("N80").Value = -5788.8572
$ws_GSM.Range("H83").Value = 3602.2917
$ws_GSM.Range("J83").Value = 3792.8572
$ws_GSM.Range("L83").Value = 18964.286
$ws_GSM.Range("N83").Value = -28948.286
$ws_GSM.Range("H107").Value = 622.75
$ws_GSM.Range("I107").Value = 519.8
$ws_GSM.Range("K107").Value = 519.8
$ws_GSM.Range("M107").Value = 1400.2
$ws_GSM.Range("H123").Value = 6404.643
$ws_GSM.Range("I123").Value = 4000
$ws_GSM.Range("J123").Value = 8488.666999999999
$ws_GSM.Range("K123").Value = 4000
$ws_GSM.Range("L123").Value = 8488.666999999999
$ws_GSM.Range("M123").Value = -1550
$ws_GSM.Range("N123").Value = -13388.667
$ws_GSM.Range("H132").Value = 22997.12
$ws_GSM.Range("I132").Value = 2514.1875
$ws_GSM.Range("J132").Value = 59411.223
$ws_GSM.Range("K132").Value = 7542.5625
$ws_GSM.Range("L132").Value = 178233.669
$ws_GSM.Range("M132").Value = -5012.5625
$ws_GSM.Range("N132").Value = -183293.669

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H61").Value = 3518.2778
$ws_LTW.Range("I61").Value = 1770.5625
$ws_LTW.Range("J61").Value = 17500
$ws_LTW.Range("K61").Value = 1770.5625
$ws_LTW.Range("L61").Value = 17500
$ws_LTW.Range("M61").Value = -1568.5625
$ws_LTW.Range("N61").Value = -17904
$ws_LTW.Range("H113").Value = 3518.2778
$ws_LTW.Range("I113").Value = 1770.5625
$ws_LTW.Range("J113").Value = 17500
$ws_LTW.Range("K113").Value = 1770.5625
$ws_LTW.Range("L113").Value = 17500
$ws_LTW.Range("M113").Value = 399.4375
$ws_LTW.Range("N113").Value = -21840
$ws_LTW.Range("H122").Value = 563678
$ws_LTW.Range("I122").Value = 820073.4399999999
$ws_LTW.Range("J122").Value = 4269.8184
$ws_LTW.Range("K122").Value = 2460220.32
$ws_LTW.Range("L122").Value = 12809.4552
$ws_LTW.Range("M122").Value = -2457770.32
$ws_LTW.Range("N122").Value = -17709.4552
$ws_LTW.Range("H130").Value = 19200
$ws_LTW.Range("J130").Value = 19200
$ws_LTW.Range("L130").Value = 19200
$ws_LTW.Range("N130").Value = -29240
$ws_LTW.Range("H136").Value = 1508.2
$ws_LTW.Range("I136").Value = 1487.7084
$ws_LTW.Range("J136").Value = 2000
$ws_LTW.Range("K136").Value = 4463.1252
$ws_LTW.Range("L136").Value = 6000
$ws_LTW.Range("M136").Value = -1913.1252
$ws_LTW.Range("N136").Value = -11100

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H46").Value = 40000
$ws_WVR.Range("J46").Value = 40000
$ws_WVR.Range("L46").Value = 40000
$ws_WVR.Range("N46").Value = -40462
$ws_WVR.Range("H107").Value = 66666868
$ws_WVR.Range("I107").Value = 83333530
$ws_WVR.Range("J107").Value = 223.33333
$ws_WVR.Range("K107").Value = 250000590
$ws_WVR.Range("L107").Value = 669.99999
$ws_WVR.Range("M107").Value = -249998670
$ws_WVR.Range("N107").Value = -4509.99999
$ws_WVR.Range("H134").Value = 40000
$ws_WVR.Range("J134").Value = 40000
$ws_WVR.Range("L134").Value = 120000
$ws_WVR.Range("N134").Value = -125070
$ws_WVR.Range("H136").Value = 29496368
$ws_WVR.Range("I136").Value = 36867748
$ws_WVR.Range("K136").Value = 110603244
$ws_WVR.Range("M136").Value = -110600694

